$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E7").Value = $true
